# member payment verification page object class
#
# Adds the locator rows used by the "MemberPayment" verification page
# object (transaction type / to / amount / description elements, plus the
# back button) to the MemberPayment worksheet, and updates the sheet's
# selection state to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MemberPayment")
$ws.Activate()

# Main_Locator_Value column first (matches how the locators were authored).
$ws.Range("C14").Value = '//*[@id="tdContents"]/form/table[1]/tbody/tr[2]/td/table/tbody/tr[2]/td[2]'
$ws.Range("C15").Value = '//*[@id="tdContents"]/form/table[1]/tbody/tr[2]/td/table/tbody/tr[3]/td[2]'
$ws.Range("C16").Value = '//*[@id="tdContents"]/form/table[1]/tbody/tr[2]/td/table/tbody/tr[4]/td[2]'
$ws.Range("C17").Value = '//*[@id="tdContents"]/form/table[1]/tbody/tr[2]/td/table/tbody/tr[5]/td[2]'

# Element_Name column.
$ws.Range("A14").Value = "ELM_TransactionType"
$ws.Range("A15").Value = "ELM_To"
$ws.Range("A16").Value = "ELM_TransactionAmount"
$ws.Range("A17").Value = "ELM_TransactionDescription"

# Main_Locator_Name column.
$ws.Range("B14").Value = "xpath"
$ws.Range("B15").Value = "xpath"
$ws.Range("B16").Value = "xpath"
$ws.Range("B17").Value = "xpath"

# New row for the back button.
$ws.Range("A18").Value = "BTN_Back"
$ws.Range("B18").Value = "xpath"
$ws.Range("C18").Value = '//*[@id="backButton"]'

# Update the worksheet selection to reflect where the new data was entered.
$ws.Range("A17").Select()
